# Update cryptos list values (Price and Volume(1h) columns) per latest snapshot.
# Cells whose new Price looks like a plain decimal number (e.g. "369.53")
# are explicitly formatted as text first so Excel keeps them as strings
# (matching the source data, which stores all Price/Volume cells as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.130.02'
$ws.Range("E2").Value = '  -1.03%  '
$ws.Range("D3").Value = '2.913.79'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '369.53'
$ws.Range("E5").Value = '  +5.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.42'
$ws.Range("E6").Value = '  -2.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.539'
$ws.Range("E7").Value = '  -2.57%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.589'
$ws.Range("E9").Value = '  -3.14%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.71'
$ws.Range("E10").Value = '  -2.32%  '
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0834'
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.39'
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").Value = '3.373.13'
$ws.Range("E14").Value = '  -0.35%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.39'
$ws.Range("E15").Value = '  -3.16%  '
$ws.Range("D16").Value = '2.916.79'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.943'
$ws.Range("E17").Value = '  -2.00%  '
$ws.Range("D18").Value = '51.030.80'
$ws.Range("E18").Value = '  -1.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.26'
$ws.Range("E19").Value = '  -4.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.21'
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.77'
$ws.Range("E21").Value = '  -4.25%  '
$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.32'
$ws.Range("E23").Value = '  -0.84%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.77'
$ws.Range("E24").Value = '  -0.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.67'
$ws.Range("E25").Value = '  -1.70%  '
$ws.Range("E26").Value = '  +1.29%  '
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '25.74'
$ws.Range("E28").Value = '  -2.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.98'
$ws.Range("E29").Value = '  -5.35%  '
$ws.Range("E30").Value = '  -0.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '9.90'
$ws.Range("E31").Value = '  -2.89%  '
$ws.Range("E32").Value = '  +3.21%  '
$ws.Range("E33").Value = '  -1.39%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '34.55'
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.80'
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("E36").Value = '  +0.63%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0422'
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.03'
$ws.Range("E38").Value = '  -3.06%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.67'
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.07'
$ws.Range("E40").Value = '  -3.25%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.85'
$ws.Range("E41").Value = '  -5.33%  '
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '22.16'
$ws.Range("E43").Value = '  -3.26%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '119.67'
$ws.Range("E44").Value = '  -0.22%  '
$ws.Range("E45").Value = '  -3.46%  '
$ws.Range("D46").Value = '2.021.55'
$ws.Range("E46").Value = '  -3.63%  '
$ws.Range("E47").Value = '  -5.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.16'
$ws.Range("E48").Value = '  -4.30%  '
$ws.Range("D49").Value = '3.196.49'
$ws.Range("E49").Value = '  -0.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.238'
$ws.Range("E50").Value = '  +0.09%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0311'
$ws.Range("E51").Value = '  -7.79%  '
